$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "25.024.85"
$ws.Range("E2").Value = "  -3.70%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.644.48"
$ws.Range("E3").Value = "  -5.56%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.51"
$ws.Range("E5").Value = "  -5.44%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4774"
$ws.Range("E7").Value = "  -5.35%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2600"
$ws.Range("E8").Value = "  -5.08%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06112"
$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07054"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.646.67"
$ws.Range("E11").Value = "  -5.49%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.68"
$ws.Range("E12").Value = "  -2.74%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5955"
$ws.Range("E13").Value = "  -8.84%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.363"
$ws.Range("E14").Value = "  -6.72%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.78"
$ws.Range("E15").Value = "  -4.87%  "

$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "25.032.72"
$ws.Range("E18").Value = "  -3.75%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.29"
$ws.Range("E19").Value = "  -5.17%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006599"
$ws.Range("E20").Value = "  -3.76%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "1.855.64"
$ws.Range("E21").Value = "  -6.10%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.336"
$ws.Range("E22").Value = "  -3.06%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.592"
$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.240"
$ws.Range("E24").Value = "  -2.39%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.67"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.93"
$ws.Range("E26").Value = "  -2.11%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.389"
$ws.Range("E27").Value = "  -7.52%  "

$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "103.86"
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.650"
$ws.Range("E29").Value = "  -7.37%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.886"
$ws.Range("E30").Value = "  -0.69%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07691"
$ws.Range("E31").Value = "  -5.74%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.568"
$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("B33").Value = "Frax"
$ws.Range("C33").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9989"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04298"
$ws.Range("E34").Value = "  -8.20%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.570"
$ws.Range("E35").Value = "  -3.27%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5940"
$ws.Range("E36").Value = "  -3.09%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9298"
$ws.Range("E37").Value = "  -6.67%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.580"
$ws.Range("E38").Value = "  -6.49%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8649"
$ws.Range("E39").Value = "  +13.42%  "

$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9993"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01513"
$ws.Range("E41").Value = "  -6.66%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.75"
$ws.Range("E42").Value = "  -1.95%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.763"
$ws.Range("E43").Value = "  -8.49%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3709"
$ws.Range("E44").Value = "  -5.22%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.682"
$ws.Range("E45").Value = "  -6.54%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1101"
$ws.Range("E46").Value = "  -4.92%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.091"
$ws.Range("E47").Value = "  -3.61%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05207"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.13"
$ws.Range("E49").Value = "  -5.06%  "

$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9999"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9989"
$ws.Range("E51").Value = "  +0.02%  "
